$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 62: S0378437118303212
$ws.Range("B62").Value2 = 'Extractive multi-document summarization using multilayer networks'
$ws.Range("A62").Value2 = 'S0378437118303212'
$ws.Range("C62").Value2 = 'Huge volumes of textual information has been produced every single day. In order to organize and understand such large datasets, in recent years, summarization techniques have become popular. These techniques aims at finding relevant, concise and non-redundant content from such a big data. While network methods have been adopted to model texts in some scenarios, a systematic evaluation of multilayer network models in the multi-document summarization task has been limited to a few studies. Here, we evaluate the performance of a multilayer-based method to select the most relevant sentences in the context of an extractive multi document summarization (MDS) task. In the adopted model, nodes represent sentences and edges are created based on the number of shared words between sentences. Differently from previous studies in multi-document summarization, we make a distinction between edges linking sentences from different documents (inter-layer) and those connecting sentences from the same document (intra-layer). As a proof of principle, our results reveal that such a discrimination between intra- and inter-layer in a multilayered representation is able to improve the quality of the generated summaries. This piece of information could be used to improve current statistical methods and related textual models.'
$ws.Range("D62").Value2 = 'Complex networks; Multilayer networks; Structure and dynamics; PageRank; Text analysis; Text summarization'
$ws.Rows.Item(62).RowHeight = 240

# Row 63: S0167739X17328431
$ws.Range("B63").Value2 = 'Analyzing scientific context of researchers and communities by using complex network and semantic technologies'
$ws.Range("A63").Value2 = 'S0167739X17328431'
$ws.Range("C63").Value2 = 'Social network communities are composed of people with common interests who influence or are influenced by themselves. In the scientific context, Scientific Social Networks are characterized as social networks that represent the social relations established by researchers. Identifying and exploring these relationships are fundamental activities to support scientific experiments. In this study, we aim to discuss the use of complex networks combined with semantic analysis in a network of scientific publications called DBLP. DBLP can be classified as big data, and its use for the analysis of connections and influences among researchers can be considered a context-aware approach. Therefore, in the present study, concepts of complex network analysis are applied to verify the level of influence among researchers, by analyzing the structure of the scientific social network under study and its communities. A bidirectional graph-based model was proposed in order to evaluate the influence of researchers, in addition to algorithms to analyze the network structure and identify scientific communities, using ontological terms and rules, considering the scientific context, and identifying new connections to promote scientific collaboration. For the identification of scientific communities, we proposed an overlapping community detection algorithm, named NetSCAN. A large scientific database (DBLP) together with digital libraries were used to evaluate the model and the algorithms in a historical research experiment. The results point to the viability and effectiveness of the proposed solution.'
$ws.Range("D63").Value2 = 'Scientific social network analysis; Semantic analysis; Overlapping community detection; Clustering algorithm; Ontology; Scientific context'
$ws.Rows.Item(63).RowHeight = 270

# Row 64: S2212827119309345
$ws.Range("B64").Value2 = 'Examination of design for large and complex network projects'
$ws.Range("A64").Value2 = 'S2212827119309345'
$ws.Range("C64").Value2 = 'This paper considers identifying critical errors problem in a complex project design with a large number of performers. A multi-stage scheme for extracting semantic data from a large number of documents with interdependent sections is offered. At the first stage, the complexity of the primary analysis of data from these sections is estimated. Then this complexity in view of the identified data from adjacent sections is clarified. Based on this, the final breakdown of a large document data set is formed into clusters with alignment of their complexity. The selection and distribution of clusters between experts are in such a way as to ensure the maximum criterion of analysis efficiency in accordance with the existing time and financial resources. The proposed approach has been applied when developing large-scale transport projects.'
$ws.Range("D64").Value2 = 'Selective processing; big data; data mining; alignment; cluster; network; heuristics; fond allocation'
$ws.Rows.Item(64).RowHeight = 150

# Row 65: S002002551930043X
$ws.Range("B65").Value2 = 'Linguistic data mining with complex networks: A stylometric-oriented approach'
$ws.Range("A65").Value2 = 'S002002551930043X'
$ws.Range("C65").Value2 = 'By representing a text by a set of words and their co-occurrences, one obtains a word-adjacency network being a reduced representation of a given language sample. In this paper, the possibility of using network representation to extract information about individual language styles of literary texts is studied. By determining selected quantitative characteristics of the networks and applying machine learning algorithms, it is possible to distinguish between texts of different authors. Within the studied set of texts, English and Polish, a properly rescaled weighted clustering coefficients and weighted degrees of only a few nodes in the word-adjacency networks are sufficient to obtain the authorship attribution accuracy over 90%. A correspondence between the text authorship and the word-adjacency network structure can therefore be found. The network representation allows to distinguish individual language styles by comparing the way the authors use particular words and punctuation marks. The presented approach can be viewed as a generalization of the authorship attribution methods based on simple lexical features. Additionally, other network parameters are studied, both local and global ones, for both the unweighted and weighted networks. Their potential to capture the writing style diversity is discussed; some differences between languages are observed.'
$ws.Range("D65").Value2 = 'Complex networks; Natural language; Data mining; Stylometry; Authorship attribution'
$ws.Rows.Item(65).RowHeight = 240

# Row 66: S1877042811024037
$ws.Range("B66").Value2 = 'KBDeX: A Platform for Exploring Discourse in Collaborative Learning'
$ws.Range("A66").Value2 = 'S1877042811024037'
$ws.Range("C66").Value2 = 'Knowledge building as defined in this study is emergent collaborative learning on ill-structured tasks. Although discourses in collaborative learning have been analyzed with traditional qualitative approaches in the learning sciences field, it is difficult to capture the group dynamics. Hence, we are trying to establish a methodology for discourse analysis in collaborative learning from the perspective of complex network science. In order to conduct this study effectively, we are currently developing an application platform, called Knowledge Building Discourse Explorer (KBDeX). The goal of this project is not only to facilitate productive communication between researchers who are concerned with research on knowledge building or emergent collaborative learning, but also to encourage students to explore their own group dynamics by themselves. KBDeX is an analysis platform to visualize network structures of discourse based on the bipartite graph of words × discourse units. KBDeX can visualize them into three different network structures of: (1) students, (2) discourse units, and (3) selected words. The users can explore these three networks with its coefficients and analyze the discourse across phases or in a and stepwise way. Using discourse which has been already analyzed with a traditional qualitative approach, we will demonstrate the beneficial attributes of the KBDeX platform.'
$ws.Range("D66").Value2 = 'Collaborative Learning; Discourse Analysis; Network Anasysis ;Bipartited Network; Application Platform'
$ws.Rows.Item(66).RowHeight = 240

# Row 67: S0378437115006792
$ws.Range("B67").Value2 = 'A measure of the diversity of concepts addressed in semantic networks'
$ws.Range("A67").Value2 = 'S0378437115006792'
$ws.Range("C67").Value2 = 'In this paper, we studied density effects in semantic networks constructed from a database of titles of papers published in scientific journals as a parameter to indicate the diversity of concepts in a journal. The proposed method essentially consists of fixing the number of titles for all of the studied scientific journals and analyzing the behavior of the density variation curves with regard to the inclusion of cliques (that is, complete networks associated with the titles). We observed that density behaves as a critically self-organized object when titles (cliques) are included in the network.'
$ws.Range("D67").Value2 = 'Complex networks; Self-organized criticality; Euclidean distances'
$ws.Rows.Item(67).RowHeight = 105

# Row 68: S1568494611004492
$ws.Range("B68").Value2 = 'iDetect: Content Based Monitoring of Complex Networks using Mobile Agents'
$ws.Range("A68").Value2 = 'S1568494611004492'
$ws.Range("C68").Value2 = 'With the evolution in computer networks over the last decade, researchers are trying to come up with efficient approaches which can help network administrator in implementing the acceptable use policy for large complex networks. In this paper we have modified An Agent Based System for Activity Monitoring on Network – ABSAMN architecture and proposed iDetect: Content Based Monitoring of Complex Network using Mobile Agents which uses the content (i.e. text, image and video) of the application for categorization purpose. iDetect is implemented in Java using Java Agent Development (JADE) framework and supports platform independence; however, the framework has been tested only on Microsoft Windows (any version) environment. We have evaluated iDetect and ABSAMN on same configuration concurrently at the university campus having four labs equipped with 60–120 number of PCs in various labs; experimental results shows that iDetect efficiently detects known/unknown illegal activities (applications/websites) running on the network as compared to ABSAMN.'
$ws.Range("D68").Value2 = 'Complex network monitoring; Content based monitoring; Mobile agent; Distributed proxy server; Collaborative Multi-Agent System'
$ws.Rows.Item(68).RowHeight = 195

# Row 69: S0020025516310283
$ws.Range("B69").Value2 = 'In narrative texts punctuation marks obey the same statistics as words'
$ws.Range("A69").Value2 = 'S0020025516310283'
$ws.Range("C69").Value2 = 'From a grammar point of view, the role of punctuation marks in a sentence is formally defined and well understood. In semantic analysis punctuation plays also a crucial role as a method of avoiding ambiguity of the meaning. A different situation can be observed in the statistical analyses of language samples, where the decision on whether the punctuation marks should be considered or should be neglected is seen rather as arbitrary and at present it belongs to a researcher’s preference. An objective of this work is to shed some light onto this problem by providing us with an answer to the question whether the punctuation marks may be treated as ordinary words and whether they should be included in any analysis of the word co-occurrences. We already know from our previous study (S. Drożdż et al., Inf. Sci. 331 (2016) 32-44) that full stops that determine the length of sentences are the main carrier of long-range correlations. Now we extend that study and analyse statistical properties of the most common punctuation marks in a few Indo-European languages, investigate their frequencies, and locate them accordingly in the Zipf rank-frequency plots as well as study their role in the word-adjacency networks. We show that, from a statistical viewpoint, the punctuation marks reveal properties that are qualitatively similar to the properties of the most frequent words like articles, conjunctions, pronouns, and prepositions. This refers to both the Zipfian analysis and the network analysis. By adding the punctuation marks to the Zipf plots, we also show that these plots that are normally described by the Zipf–Mandelbrot distribution largely restore the power-law Zipfian behaviour for the most frequent items. Our results indicate that the punctuation marks can fruitfully be considered in the linguistic studies as their inclusion effectively extends dimensionality of an analysis and, therefore, it opens more space for possible manifestation of some previously unobserved effects.'
$ws.Range("D69").Value2 = 'Punctuation; Word-adjacency networks; Complex networks; Word-frequency distribution'
$ws.Rows.Item(69).RowHeight = 345

# Row 70: S0378437106006881
$ws.Range("B70").Value2 = 'Strong correlations between text quality and complex networks features'
$ws.Range("A70").Value2 = 'S0378437106006881'
$ws.Range("C70").Value2 = 'Concepts of complex networks have been used to obtain metrics that were correlated to text quality established by scores assigned by human judges. Texts produced by high-school students in Portuguese were represented as scale-free networks (word adjacency model), from which typical network features such as the in/outdegree, clustering coefficient and shortest path were obtained. Another metric was derived from the dynamics of the network growth, based on the variation of the number of connected components. The scores assigned by the human judges according to three text quality criteria (coherence and cohesion, adherence to standard writing conventions and theme adequacy/development) were correlated with the network measurements. Text quality for all three criteria was found to decrease with increasing average values of outdegrees, clustering coefficient and deviation from the dynamics of network growth. Among the criteria employed, cohesion and coherence showed the strongest correlation, which probably indicates that the network measurements are able to capture how the text is developed in terms of the concepts represented by the nodes in the networks. Though based on a particular set of texts and specific language, the results presented here point to potential applications in other instances of text analysis.'
$ws.Range("D70").Value2 = 'Complex networks; Text analysis; Networks measurements'
$ws.Rows.Item(70).RowHeight = 240

# Row 71: S0378437118308343
$ws.Range("B71").Value2 = 'An image analysis approach to text analytics based on complex networks'
$ws.Range("A71").Value2 = 'S0378437118308343'
$ws.Range("C71").Value2 = 'Text network analysis has received increasing attention as a consequence of its wide range of applications. In this study, we extend a previous work founded on the study of topological features of mesoscopic networks. Here, the geometrical properties of visualized networks are quantified by using several image analysis techniques. Such properties are used to probe the networks characteristics in terms of authorship. It was found that the visual features account for performance similar to that achieved by using topological measurements. Also, we combined and compared the two types of features, topological and geometrical, and the results suggest that the information provided by network topology and image features are complementary.'
$ws.Range("D71").Value2 = 'Complex networks; Mesoscopic networks; Authorship; Text mining; Computer vision'
$ws.Rows.Item(71).RowHeight = 135

# Update selection to the new last cell, matching the final saved UI state
[void]$ws.Range("D71").Select()

